$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45204 = 2023-10-05) for every
# data row (2 through 303). Bump it by one day to 45205 (2023-10-06) across the
# whole column, matching the diff.
for ($row = 2; $row -le 303; $row++) {
    $ws.Cells.Item($row, 3).Value = 45205
}
